$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.045368349355385
$ws.Range("D2").Value = 1.052048691316282
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.059592022044555
$ws.Range("I2").Value = 1.043169153626973
$ws.Range("J2").Value = 1.050428456872102
$ws.Range("K2").Value = 1.054798501661742
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.062321106454126
$ws.Range("N2").Value = 1.02064360924143

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.046358502311191
$ws.Range("D3").Value = 1.052851754081617
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.060577287711282
$ws.Range("I3").Value = 1.04343211435256
$ws.Range("J3").Value = 1.051066042704576
$ws.Range("K3").Value = 1.055414357420766
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.063120199222103
$ws.Range("N3").Value = 1.020858664021513

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046999347063445
$ws.Range("D4").Value = 1.053371532121558
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.061215396576055
$ws.Range("I4").Value = 1.043601173553645
$ws.Range("J4").Value = 1.051478138309381
$ws.Range("K4").Value = 1.055812345466
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.063637223167919
$ws.Range("N4").Value = 1.020997582458498

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.047268793587151
$ws.Range("D5").Value = 1.053590080078998
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.06148379461951
$ws.Range("I5").Value = 1.043671983804806
$ws.Range("J5").Value = 1.05165127140101
$ws.Range("K5").Value = 1.055979536481087
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.063854569130765
$ws.Range("N5").Value = 1.02105592691647

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.047314036907266
$ws.Range("D6").Value = 1.053626777167469
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.061528867845834
$ws.Range("I6").Value = 1.04368385777007
$ws.Range("J6").Value = 1.051680334630659
$ws.Range("K6").Value = 1.056007601347789
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.063891061813776
$ws.Range("N6").Value = 1.021065719866201

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.047002947282537
$ws.Range("D7").Value = 1.053374452240461
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.061218982386278
$ws.Range("I7").Value = 1.043602120754516
$ws.Range("J7").Value = 1.051480452162619
$ws.Range("K7").Value = 1.055814579965242
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.063640127399063
$ws.Range("N7").Value = 1.02099836228372

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045702944957851
$ws.Range("D8").Value = 1.05232006004288
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.0599248773989
$ws.Range("I8").Value = 1.043258248746006
$ws.Range("J8").Value = 1.050644027944677
$ws.Range("K8").Value = 1.055006738593913
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.0625911719607
$ws.Range("N8").Value = 1.02071633679934

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043413338651473
$ws.Range("D9").Value = 1.050463218992363
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.057648948229222
$ws.Range("I9").Value = 1.04264393866365
$ws.Range("J9").Value = 1.049166604473236
$ws.Range("K9").Value = 1.0535793266671
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.060742484803834
$ws.Range("N9").Value = 1.020217572898167

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041887739302098
$ws.Range("D10").Value = 1.049226138440306
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.05613470152473
$ws.Range("I10").Value = 1.042228794609955
$ws.Range("J10").Value = 1.048179309153318
$ws.Range("K10").Value = 1.052625137737065
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.05950987069093
$ws.Range("N10").Value = 1.019883869085494

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041227332786072
$ws.Range("D11").Value = 1.048690672331202
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.055479746444942
$ws.Range("I11").Value = 1.042047707254182
$ws.Range("J11").Value = 1.047751249087295
$ws.Range("K11").Value = 1.052211358349703
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.058976106380616
$ws.Range("N11").Value = 1.019739091317705

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040982056879285
$ws.Range("D12").Value = 1.048491806901958
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.055236576402252
$ws.Range("I12").Value = 1.041980244204204
$ws.Range("J12").Value = 1.047592165621366
$ws.Range("K12").Value = 1.052057571236096
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.058777837996938
$ws.Range("N12").Value = 1.019685272288327

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041034668101606
$ws.Range("D13").Value = 1.048534462817612
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.055288732238373
$ws.Range("I13").Value = 1.041994724264888
$ws.Range("J13").Value = 1.047626293309127
$ws.Range("K13").Value = 1.052090563204871
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.058820367434741
$ws.Range("N13").Value = 1.01969681855533

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041207057620861
$ws.Range("D14").Value = 1.04867423342552
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.055459643696829
$ws.Range("I14").Value = 1.042042134798375
$ws.Range("J14").Value = 1.047738100886943
$ws.Range("K14").Value = 1.052198648112846
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.058959717543604
$ws.Range("N14").Value = 1.019734643477919

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04131327627357
$ws.Range("D15").Value = 1.048760354778536
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.05556496240953
$ws.Range("I15").Value = 1.042071319612213
$ws.Range("J15").Value = 1.047806978254027
$ws.Range("K15").Value = 1.052265230745475
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.059045575159322
$ws.Range("N15").Value = 1.01975794308633

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041931572288888
$ws.Range("D16").Value = 1.049261679799953
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.056178184044563
$ws.Range("I16").Value = 1.042240784840051
$ws.Range("J16").Value = 1.048207706440885
$ws.Range("K16").Value = 1.052652586119055
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.059545294199121
$ws.Range("N16").Value = 1.019893471590337

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042319463551775
$ws.Range("D17").Value = 1.049576201289909
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.056563036035462
$ws.Range("I17").Value = 1.042346730743013
$ws.Range("J17").Value = 1.048458924449948
$ws.Range("K17").Value = 1.052895401149961
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.059858745996276
$ws.Range("N17").Value = 1.01997840974535

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042545732093715
$ws.Range("D18").Value = 1.049759675398607
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.056787583473404
$ws.Range("I18").Value = 1.042408399118713
$ws.Range("J18").Value = 1.048605402116746
$ws.Range("K18").Value = 1.053036972227908
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.060041573780298
$ws.Range("N18").Value = 1.0200279255284

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042622886887806
$ws.Range("D19").Value = 1.049822238534638
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.05686416016193
$ws.Range("I19").Value = 1.04242940470921
$ws.Range("J19").Value = 1.048655338119951
$ws.Range("K19").Value = 1.053085234314244
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.060103912744317
$ws.Range("N19").Value = 1.020044804509578

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042277844613174
$ws.Range("D20").Value = 1.04954245413499
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.056521737835934
$ws.Range("I20").Value = 1.042335376995992
$ws.Range("J20").Value = 1.048431976670555
$ws.Range("K20").Value = 1.052869355475552
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.059825115926661
$ws.Range("N20").Value = 1.019969299497578

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041156292416202
$ws.Range("D21").Value = 1.048633073616206
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.055409311505555
$ws.Range("I21").Value = 1.042028179073168
$ws.Range("J21").Value = 1.04770517862689
$ws.Range("K21").Value = 1.052166822310024
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.058918682530652
$ws.Range("N21").Value = 1.019723506141941

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040451292912901
$ws.Range("D22").Value = 1.048061486990912
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.054710518430471
$ws.Range("I22").Value = 1.041833879495807
$ws.Range("J22").Value = 1.047247732444781
$ws.Range("K22").Value = 1.051724585433815
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.058348745831917
$ws.Range("N22").Value = 1.01956872251778

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040825010784879
$ws.Range("D23").Value = 1.048364478835099
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.055080901576645
$ws.Range("I23").Value = 1.041936990491848
$ws.Range("J23").Value = 1.047490278617246
$ws.Range("K23").Value = 1.051959073293739
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.058650882329327
$ws.Range("N23").Value = 1.0196507992615

# Row 24
$ws.Range("B24").Value = 1.019999999999999
$ws.Range("C24").Value = 1.042296650363957
$ws.Range("D24").Value = 1.049557702962855
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.056540398499029
$ws.Range("I24").Value = 1.042340507661282
$ws.Range("J24").Value = 1.048444153377544
$ws.Range("K24").Value = 1.05288112457644
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.059840311917866
$ws.Range("N24").Value = 1.019973416110467

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044005116547456
$ws.Range("D25").Value = 1.050943117328685
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.058236798277751
$ws.Range("I25").Value = 1.042803741704036
$ws.Range("J25").Value = 1.049548969746077
$ws.Range("K25").Value = 1.053948804412902
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.061220445280507
$ws.Range("N25").Value = 1.020346726964031
